$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "248.05"
$cell.Style = $origStyle

$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "21.94"
$cell.Style = $origStyle

$cell = $ws.Range("D4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.331"
$cell.Style = $origStyle

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.05633"
$cell.Style = $origStyle

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.426"
$cell.Style = $origStyle

$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.371"
$cell.Style = $origStyle

$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.8136"
$cell.Style = $origStyle

$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.9347"
$cell.Style = $origStyle

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.1431"
$cell.Style = $origStyle

$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.07601"
$cell.Style = $origStyle

$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.03225"
$cell.Style = $origStyle

$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.03092"
$cell.Style = $origStyle

$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.09301"
$cell.Style = $origStyle

$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.001590"
$cell.Style = $origStyle

$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.04729"
$cell.Style = $origStyle

$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0005777"
$cell.Style = $origStyle

$ws.Range("E18").Value = "17OneONE"

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.006273"
$cell.Style = $origStyle

$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.005056"
$cell.Style = $origStyle

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.001034"
$cell.Style = $origStyle

$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0001499"
$cell.Style = $origStyle

$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.762"
$cell.Style = $origStyle

$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.146"
$cell.Style = $origStyle

$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.3306"
$cell.Style = $origStyle

$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0002999"
$cell.Style = $origStyle

$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.03953"
$cell.Style = $origStyle

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.006817"
$cell.Style = $origStyle

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.008813"
$cell.Style = $origStyle

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.00005578"
$cell.Style = $origStyle

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0005497"
$cell.Style = $origStyle

$ws.Range("E47").Value = "46ACDXExchangeACXTWorstin24h"

$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.7797"
$cell.Style = $origStyle

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.1780"
$cell.Style = $origStyle
